{"js": "// 1) Replace the employee name in the title paragraph.\nconst body = context.document.body;\nconst nameResults = body.search(\"\u041f\u0435\u0442\u0440\u043e\u0432 \u0418\u0432\u0430\u043d \u041f\u0435\u0442\u0440\u043e\u0432\u0438\u0447\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"\u0412\u0438\u043a\u0442\u043e\u0440 \u041a\u0443\u0440\u043e\u0447\u043a\u0438\u043d \u0410\u043b\u0435\u043a\u0441\u0430\u043d\u0434\u0440\u043e\u0432\u0438\u0447\", \"Replace\");\n  await context.sync();\n}\n\n// 2) The history table: remove the \"\u0418\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0439 \u043f\u043e\u0440\u0442\u0430\u043b\" row entirely, and\n//    update the date on the remaining (\"\u0427\u0430\u0442\") row from 03.04.2016 to\n//    02.05.2016.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load each row's cell text so we can find the right rows regardless of\n// ordering.\nfor (const row of rows.items) {\n  row.cells.load(\"items/body/text\");\n}\nawait context.sync();\n\nlet targetRowToDelete = null;\n\nfor (const row of rows.items) {\n  const cellTexts = row.cells.items.map((c) => c.body.text.trim());\n  if (cellTexts.some((t) => t.includes(\"\u0418\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0439 \u043f\u043e\u0440\u0442\u0430\u043b\"))) {\n    targetRowToDelete = row;\n  }\n}\n\nif (targetRowToDelete) {\n  targetRowToDelete.delete();\n  await context.sync();\n}\n\n// The row collection is stale after the structural edit above (deleting a\n// row shifts indices), so re-fetch it before locating the remaining row.\nconst rows2 = table.rows;\nrows2.load(\"items\");\nawait context.sync();\n\nfor (const row of rows2.items) {\n  row.cells.load(\"items/body/text\");\n}\nawait context.sync();\n\nlet targetRowToUpdate = null;\nfor (const row of rows2.items) {\n  const cellTexts = row.cells.items.map((c) => c.body.text.trim());\n  if (cellTexts.some((t) => t.includes(\"\u0427\u0430\u0442\"))) {\n    targetRowToUpdate = row;\n  }\n}\n\nif (targetRowToUpdate) {\n  const dateCell = targetRowToUpdate.cells.items[0];\n  const dateResults = dateCell.body.search(\"03.04.2016\", { matchCase: true });\n  dateResults.load(\"items\");\n  await context.sync();\n  if (dateResults.items.length > 0) {\n    dateResults.items[0].insertText(\"02.05.2016\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Replace the employee name in the title paragraph.\n$find = $d.Content.Find\n$find.Execute(\"\u041f\u0435\u0442\u0440\u043e\u0432 \u0418\u0432\u0430\u043d \u041f\u0435\u0442\u0440\u043e\u0432\u0438\u0447\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0412\u0438\u043a\u0442\u043e\u0440 \u041a\u0443\u0440\u043e\u0447\u043a\u0438\u043d \u0410\u043b\u0435\u043a\u0441\u0430\u043d\u0434\u0440\u043e\u0432\u0438\u0447\", 1)\n\n# 2) The history table: remove the \"\u0418\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0439 \u043f\u043e\u0440\u0442\u0430\u043b\" row entirely, and\n#    update the date on the remaining (\"\u0427\u0430\u0442\") row from 03.04.2016 to\n#    02.05.2016.\n$t = $d.Tables.Item(1)\n\n# Walk the rows bottom-up (deleting shifts indices below it but not above)\n# and look up each row's second-column text to find the row to drop.\n$rowToDelete = 0\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $cellText = $t.Cell($i, 2).Range.Text\n    if ($cellText -like \"*\u0418\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0439 \u043f\u043e\u0440\u0442\u0430\u043b*\") {\n        $rowToDelete = $i\n    }\n}\nif ($rowToDelete -gt 0) {\n    $t.Rows($rowToDelete).Delete()\n}\n\n# Re-locate the remaining data row (now holding \"\u0427\u0430\u0442\") and fix up its date.\n$rowToUpdate = 0\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $cellText = $t.Cell($i, 2).Range.Text\n    if ($cellText -like \"*\u0427\u0430\u0442*\") {\n        $rowToUpdate = $i\n    }\n}\nif ($rowToUpdate -gt 0) {\n    $dateRange = $t.Cell($rowToUpdate, 1).Range\n    $dateFind = $dateRange.Find\n    $dateFind.Execute(\"03.04.2016\", $false, $false, $false, $false, $false, $true, 1, $false, \"02.05.2016\", 1)\n}\n"}
